$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.953.43'
$ws.Range("E2").Value = '  +0.45%  '

# Row 3
$ws.Range("D3").Value = '1.640.38'
$ws.Range("E3").Value = '  +0.32%  '

# Row 4
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '214.69'

# Row 6
$ws.Range("D6").Value = '0.5098'
$ws.Range("E6").Value = '  +1.58%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("E8").Value = '  +0.22%  '

# Row 9
$ws.Range("D9").Value = '0.06366'
$ws.Range("E9").Value = '  +0.12%  '

# Row 10
$ws.Range("D10").Value = '19.50'
$ws.Range("E10").Value = '  +0.77%  '

# Row 11
$ws.Range("D11").Value = '0.07765'
$ws.Range("E11").Value = '  -0.20%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.657.23'
$ws.Range("E12").Value = '  +0.99%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.280'
$ws.Range("E13").Value = '  +0.99%  '

# Row 14
$ws.Range("D14").Value = '0.5440'
$ws.Range("E14").Value = '  +0.79%  '

# Row 15
$ws.Range("D15").Value = '0.0₅7740'
$ws.Range("E15").Value = '  -1.49%  '

# Row 16
$ws.Range("D16").Value = '64.25'
$ws.Range("E16").Value = '  +0.15%  '

# Row 17
$ws.Range("D17").Value = '25.986.87'
$ws.Range("E17").Value = '  +0.45%  '

# Row 18
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.12%  '

# Row 19
$ws.Range("D19").Value = '195.93'
$ws.Range("E19").Value = '  +0.49%  '

# Row 20
$ws.Range("D20").Value = '4.422'
$ws.Range("E20").Value = '  +1.60%  '

# Row 21
$ws.Range("E21").Value = '  +0.67%  '

# Row 22
$ws.Range("D22").Value = '6.048'
$ws.Range("E22").Value = '  +1.81%  '

# Row 23
$ws.Range("E23").Value = '  -0.05%  '

# Row 24
$ws.Range("D24").Value = '1.882'
$ws.Range("E24").Value = '  -0.69%  '

# Row 25
$ws.Range("D25").Value = '141.18'
$ws.Range("E25").Value = '  +1.11%  '

# Row 26
$ws.Range("D26").Value = '0.1192'
$ws.Range("E26").Value = '  +5.66%  '

# Row 27
$ws.Range("D27").Value = '6.849'
$ws.Range("E27").Value = '  +0.98%  '

# Row 28
$ws.Range("D28").Value = '15.62'
$ws.Range("E28").Value = '  +0.11%  '

# Row 29
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("D30").Value = '0.04917'
$ws.Range("E30").Value = '  +1.46%  '

# Row 31
$ws.Range("D31").Value = '3.250'
$ws.Range("E31").Value = '  +0.54%  '

# Row 32
$ws.Range("D32").Value = '3.173'
$ws.Range("E32").Value = '  +0.45%  '

# Row 33
$ws.Range("D33").Value = '1.528'
$ws.Range("E33").Value = '  +0.17%  '

# Row 34
$ws.Range("D34").Value = '2.370'
$ws.Range("E34").Value = '  +0.38%  '

# Row 35
$ws.Range("D35").Value = '0.8920'
$ws.Range("E35").Value = '  +1.13%  '

# Row 36
$ws.Range("D36").Value = '1.150.60'
$ws.Range("E36").Value = '  +2.32%  '

# Row 37
$ws.Range("D37").Value = '2.585'
$ws.Range("E37").Value = '  -0.43%  '

# Row 38
$ws.Range("D38").Value = '0.5425'
$ws.Range("E38").Value = '  -1.33%  '

# Row 39
$ws.Range("E39").Value = '  +0.01%  '

# Row 40
$ws.Range("D40").Value = '1.002'
$ws.Range("E40").Value = '  -0.04%  '

# Row 41
$ws.Range("D41").Value = '2.528'
$ws.Range("E41").Value = '  -0.90%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.8111'
$ws.Range("E42").Value = '  +0.12%  '

# Row 43
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D43").Value = '0.0₈126'
$ws.Range("E43").Value = '  +3.23%  '

# Row 44
$ws.Range("D44").Value = '99.01'
$ws.Range("E44").Value = '  -0.17%  '

# Row 45
$ws.Range("D45").Value = '5.445'
$ws.Range("E45").Value = '  -3.49%  '

# Row 46
$ws.Range("D46").Value = '1.778.20'
$ws.Range("E46").Value = '  +0.26%  '

# Row 47
$ws.Range("D47").Value = '0.4528'
$ws.Range("E47").Value = '  +0.31%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '54.88'
$ws.Range("E48").Value = '  -0.09%  '

# Row 49
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '0.9991'
$ws.Range("E49").Value = '  -0.34%  '

# Row 50
$ws.Range("D50").Value = '0.05053'
$ws.Range("E50").Value = '  +0.40%  '

# Row 51
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  -0.26%  '

